$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 496.91306
$ws.Range("I2").Value = 67.36842
$ws.Range("K2").Value = 67.36842
$ws.Range("M2").Value = 45.63158
$ws.Range("H33").Value = 2107103.5
$ws.Range("I33").Value = 2668613
$ws.Range("K33").Value = 2668613
$ws.Range("M33").Value = -2668384
$ws.Range("H41").Value = 23812786
$ws.Range("I41").Value = 828.4286
$ws.Range("J41").Value = 35718764
$ws.Range("K41").Value = 828.4286
$ws.Range("L41").Value = 35718764
$ws.Range("M41").Value = -388.4286
$ws.Range("N41").Value = -35719644
$ws.Range("H62").Value = 5438709
$ws.Range("I62").Value = 7815233
$ws.Range("K62").Value = 7815233
$ws.Range("M62").Value = -7814609
$ws.Range("H63").Value = 52333
$ws.Range("J63").Value = 52333
$ws.Range("L63").Value = 52333
$ws.Range("N63").Value = -53581
$ws.Range("H64").Value = 6631.643
$ws.Range("J64").Value = 8105.5
$ws.Range("L64").Value = 8105.5
$ws.Range("N64").Value = -8601.5
$ws.Range("H65").Value = 5438709
$ws.Range("I65").Value = 7815233
$ws.Range("K65").Value = 39076165
$ws.Range("M65").Value = -39073045
$ws.Range("H66").Value = 52333
$ws.Range("J66").Value = 52333
$ws.Range("L66").Value = 156999
$ws.Range("N66").Value = -163239
$ws.Range("H67").Value = 6631.643
$ws.Range("J67").Value = 8105.5
$ws.Range("L67").Value = 8105.5
$ws.Range("N67").Value = -9821.5
$ws.Range("H70").Value = 250874.75
$ws.Range("I70").Value = 500
$ws.Range("J70").Value = 334333
$ws.Range("K70").Value = 1500
$ws.Range("L70").Value = 1002999
$ws.Range("M70").Value = -1230
$ws.Range("N70").Value = -1003539
$ws.Range("H73").Value = 250874.75
$ws.Range("I73").Value = 500
$ws.Range("J73").Value = 334333
$ws.Range("K73").Value = 1500
$ws.Range("L73").Value = 1002999
$ws.Range("M73").Value = -564
$ws.Range("N73").Value = -1004871
$ws.Range("H86").Value = 3764005.2
$ws.Range("I86").Value = 3749.8333
$ws.Range("K86").Value = 3749.8333
$ws.Range("M86").Value = -2626.8333
$ws.Range("H89").Value = 3764005.2
$ws.Range("I89").Value = 3749.8333
$ws.Range("K89").Value = 18749.1665
$ws.Range("M89").Value = -13133.1665
$ws.Range("H92").Value = 146.13333
$ws.Range("I92").Value = 102.545456
$ws.Range("J92").Value = 266
$ws.Range("K92").Value = 102.545456
$ws.Range("L92").Value = 266
$ws.Range("M92").Value = 1145.454544
$ws.Range("N92").Value = -2762
$ws.Range("H97").Value = 3180.6667
$ws.Range("J97").Value = 3946.8572
$ws.Range("L97").Value = 11840.5716
$ws.Range("N97").Value = -12832.5716
$ws.Range("H98").Value = 2349.3103
$ws.Range("I98").Value = 1509.7826
$ws.Range("K98").Value = 1509.7826
$ws.Range("M98").Value = -11.7826
$ws.Range("H101").Value = 1293
$ws.Range("I101").Value = 1341.8334
$ws.Range("K101").Value = 4025.5002
$ws.Range("M101").Value = -2403.5002
$ws.Range("H103").Value = 960.25
$ws.Range("I103").Value = 867.7
$ws.Range("J103").Value = 1423
$ws.Range("K103").Value = 2603.1
$ws.Range("L103").Value = 4269
$ws.Range("M103").Value = -2017.1
$ws.Range("N103").Value = -5441
$ws.Range("H106").Value = 1849.8
$ws.Range("I106").Value = 1849.8
$ws.Range("K106").Value = 1849.8
$ws.Range("M106").Value = -1218.8
$ws.Range("H118").Value = 793.2222
$ws.Range("I118").Value = 575.7143
$ws.Range("J118").Value = 1554.5
$ws.Range("K118").Value = 1727.1429
$ws.Range("L118").Value = 4663.5
$ws.Range("M118").Value = -70.14289999999983
$ws.Range("N118").Value = -7977.5
$ws.Range("H122").Value = 2349.3103
$ws.Range("I122").Value = 1509.7826
$ws.Range("K122").Value = 4529.3478
$ws.Range("M122").Value = -2079.3478
$ws.Range("H129").Value = 8514.612999999999
$ws.Range("I129").Value = 7406
$ws.Range("K129").Value = 22218
$ws.Range("M129").Value = -17218
$ws.Range("H132").Value = 2304.16
$ws.Range("I132").Value = 1917.5652
$ws.Range("J132").Value = 6750
$ws.Range("K132").Value = 5752.6956
$ws.Range("L132").Value = 20250
$ws.Range("M132").Value = -3222.6956
$ws.Range("N132").Value = -25310
$ws.Range("H137").Value = 331488.38
$ws.Range("I137").Value = 220452.12
$ws.Range("J137").Value = 671999.5600000001
$ws.Range("K137").Value = 661356.36
$ws.Range("L137").Value = 2015998.68
$ws.Range("M137").Value = -658806.36
$ws.Range("N137").Value = -2021098.68
$ws.Range("H138").Value = 4857.5835
$ws.Range("J138").Value = 6611.6665
$ws.Range("L138").Value = 19834.9995
$ws.Range("N138").Value = -30114.9995
$ws.Range("H141").Value = 2431.1455
$ws.Range("I141").Value = 1140.8096
$ws.Range("J141").Value = 6599.923
$ws.Range("K141").Value = 3422.4288
$ws.Range("L141").Value = 19799.769
$ws.Range("M141").Value = 1757.5712
$ws.Range("N141").Value = -30159.769

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3737.1714
$ws.Range("I2").Value = 3737.1714
$ws.Range("K2").Value = 3737.1714
$ws.Range("M2").Value = -3624.1714
$ws.Range("H21").Value = 5149.5
$ws.Range("I21").Value = 3532.6667
$ws.Range("K21").Value = 3532.6667
$ws.Range("M21").Value = -3158.6667
$ws.Range("H74").Value = 1421.6444
$ws.Range("I74").Value = 1153.0488
$ws.Range("K74").Value = 1153.0488
$ws.Range("M74").Value = -279.0488
$ws.Range("H77").Value = 1421.6444
$ws.Range("I77").Value = 1153.0488
$ws.Range("K77").Value = 5765.244000000001
$ws.Range("M77").Value = -1397.244000000001
$ws.Range("H97").Value = 4419.8
$ws.Range("I97").Value = 10055
$ws.Range("K97").Value = 10055
$ws.Range("M97").Value = -9559
$ws.Range("H116").Value = 3737.1714
$ws.Range("I116").Value = 3737.1714
$ws.Range("K116").Value = 3737.1714
$ws.Range("M116").Value = -1443.1714
$ws.Range("H122").Value = 3509.6365
$ws.Range("I122").Value = 1773.1111
$ws.Range("J122").Value = 6267.647
$ws.Range("K122").Value = 5319.3333
$ws.Range("L122").Value = 18802.941
$ws.Range("M122").Value = -2869.3333
$ws.Range("N122").Value = -23702.941

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3737.1714
$ws.Range("I3").Value = 3737.1714
$ws.Range("K3").Value = 3737.1714
$ws.Range("M3").Value = -3623.1714
$ws.Range("H50").Value = 80766.336
$ws.Range("J50").Value = 80766.336
$ws.Range("L50").Value = 80766.336
$ws.Range("N50").Value = -81914.336
$ws.Range("H86").Value = 1890744.5
$ws.Range("I86").Value = 2126887.5
$ws.Range("J86").Value = 1600
$ws.Range("K86").Value = 2126887.5
$ws.Range("L86").Value = 1600
$ws.Range("M86").Value = -2125764.5
$ws.Range("N86").Value = -3846
$ws.Range("H88").Value = 30000
$ws.Range("J88").Value = 30000
$ws.Range("L88").Value = 30000
$ws.Range("N88").Value = -30812
$ws.Range("H89").Value = 1890744.5
$ws.Range("I89").Value = 2126887.5
$ws.Range("J89").Value = 1600
$ws.Range("K89").Value = 10634437.5
$ws.Range("L89").Value = 8000
$ws.Range("M89").Value = -10628821.5
$ws.Range("N89").Value = -19232
$ws.Range("H91").Value = 30000
$ws.Range("J91").Value = 30000
$ws.Range("L91").Value = 30000
$ws.Range("N91").Value = -32808
$ws.Range("H94").Value = 1261.4375
$ws.Range("I94").Value = 1071.1818
$ws.Range("J94").Value = 1680
$ws.Range("K94").Value = 1071.1818
$ws.Range("L94").Value = 1680
$ws.Range("M94").Value = -620.1818000000001
$ws.Range("N94").Value = -2582
$ws.Range("H132").Value = 50000
$ws.Range("J132").Value = 50000
$ws.Range("L132").Value = 50000
$ws.Range("N132").Value = -60120

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 57916.895
$ws.Range("I31").Value = 1578.875
$ws.Range("K31").Value = 1578.875
$ws.Range("M31").Value = -1283.875
$ws.Range("H34").Value = 57916.895
$ws.Range("I34").Value = 1578.875
$ws.Range("K34").Value = 1578.875
$ws.Range("M34").Value = -1376.875
$ws.Range("H58").Value = 124199.39
$ws.Range("I58").Value = 157881.36
$ws.Range("K58").Value = 157881.36
$ws.Range("M58").Value = -157678.36
$ws.Range("I62").Value = 3133
$ws.Range("J62").Value = 5369.75
$ws.Range("K62").Value = 3133
$ws.Range("L62").Value = 5369.75
$ws.Range("M62").Value = -2509
$ws.Range("N62").Value = -6617.75
$ws.Range("I65").Value = 3133
$ws.Range("J65").Value = 5369.75
$ws.Range("K65").Value = 15665
$ws.Range("L65").Value = 26848.75
$ws.Range("M65").Value = -12545
$ws.Range("N65").Value = -33088.75
$ws.Range("H86").Value = 80901
$ws.Range("I86").Value = 99998.5
$ws.Range("J86").Value = 68169.336
$ws.Range("K86").Value = 99998.5
$ws.Range("L86").Value = 68169.336
$ws.Range("M86").Value = -98875.5
$ws.Range("N86").Value = -70415.336
$ws.Range("H89").Value = 80901
$ws.Range("I89").Value = 99998.5
$ws.Range("J89").Value = 68169.336
$ws.Range("K89").Value = 499992.5
$ws.Range("L89").Value = 340846.68
$ws.Range("M89").Value = -494376.5
$ws.Range("N89").Value = -352078.68
$ws.Range("H107").Value = 790
$ws.Range("J107").Value = 3000
$ws.Range("L107").Value = 3000
$ws.Range("N107").Value = -6840
$ws.Range("H132").Value = 1891.8767
$ws.Range("I132").Value = 1560.4531
$ws.Range("J132").Value = 4248.6665
$ws.Range("K132").Value = 4681.3593
$ws.Range("L132").Value = 12745.9995
$ws.Range("M132").Value = -2151.3593
$ws.Range("N132").Value = -17805.9995
$ws.Range("H134").Value = 346947.28
$ws.Range("I134").Value = 219659.84
$ws.Range("K134").Value = 658979.52
$ws.Range("M134").Value = -656444.52
$ws.Range("H136").Value = 124199.39
$ws.Range("I136").Value = 157881.36
$ws.Range("K136").Value = 473644.08
$ws.Range("M136").Value = -471094.08

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 832635.0600000001
$ws.Range("J5").Value = 1431382.8
$ws.Range("L5").Value = 4294148.4
$ws.Range("N5").Value = -4294372.4
$ws.Range("H13").Value = 3175.25
$ws.Range("I13").Value = 101
$ws.Range("J13").Value = 4200
$ws.Range("K13").Value = 303
$ws.Range("L13").Value = 12600
$ws.Range("M13").Value = -135
$ws.Range("N13").Value = -12936
$ws.Range("H15").Value = 958.9091
$ws.Range("I15").Value = 1717.5
$ws.Range("J15").Value = 48.6
$ws.Range("K15").Value = 5152.5
$ws.Range("L15").Value = 145.8
$ws.Range("M15").Value = -5012.5
$ws.Range("N15").Value = -425.8
$ws.Range("H38").Value = 68.30768999999999
$ws.Range("I38").Value = 42.8
$ws.Range("J38").Value = 153.33333
$ws.Range("K38").Value = 128.4
$ws.Range("L38").Value = 459.99999
$ws.Range("M38").Value = 218.6
$ws.Range("N38").Value = -1153.99999
$ws.Range("H81").Value = 49912.23
$ws.Range("I81").Value = 17999.5
$ws.Range("J81").Value = 55714.547
$ws.Range("K81").Value = 53998.5
$ws.Range("L81").Value = 167143.641
$ws.Range("M81").Value = -52875.5
$ws.Range("N81").Value = -169389.641
$ws.Range("H84").Value = 49912.23
$ws.Range("I84").Value = 17999.5
$ws.Range("J84").Value = 55714.547
$ws.Range("K84").Value = 161995.5
$ws.Range("L84").Value = 501430.923
$ws.Range("M84").Value = -156379.5
$ws.Range("N84").Value = -512662.923
$ws.Range("H111").Value = 333875.34
$ws.Range("I111").Value = 500313
$ws.Range("K111").Value = 1500939
$ws.Range("M111").Value = -1497872
$ws.Range("H113").Value = 1950685.6
$ws.Range("J113").Value = 1489.5333
$ws.Range("L113").Value = 4468.5999
$ws.Range("N113").Value = -8808.599900000001
$ws.Range("H131").Value = 6099994
$ws.Range("J131").Value = 22786.705
$ws.Range("L131").Value = 68360.11500000001
$ws.Range("N131").Value = -78440.11500000001
$ws.Range("H135").Value = 832635.0600000001
$ws.Range("J135").Value = 1431382.8
$ws.Range("L135").Value = 12882445.2
$ws.Range("N135").Value = -12887515.2
$ws.Range("H137").Value = 4128.143
$ws.Range("J137").Value = 7766.6665
$ws.Range("L137").Value = 23299.9995
$ws.Range("N137").Value = -33499.99950000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 150000
$ws.Range("I18").Value = 150000
$ws.Range("J18").Value = 150000
$ws.Range("K18").Value = 150000
$ws.Range("L18").Value = 150000
$ws.Range("M18").Value = -149707
$ws.Range("N18").Value = -150586
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H29").Value = 20000
$ws.Range("J29").Value = 20000
$ws.Range("L29").Value = 20000
$ws.Range("N29").Value = -20580
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H102").Value = 2404.1428
$ws.Range("I102").Value = 1673.0555
$ws.Range("K102").Value = 1673.0555
$ws.Range("M102").Value = -51.05549999999994
$ws.Range("H132").Value = 585176
$ws.Range("I132").Value = 629271.5600000001
$ws.Range("J132").Value = 349999.66
$ws.Range("K132").Value = 1887814.68
$ws.Range("L132").Value = 1049998.98
$ws.Range("M132").Value = -1885284.68
$ws.Range("N132").Value = -1055058.98

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 471.57144
$ws.Range("I16").Value = 401.70587
$ws.Range("K16").Value = 401.70587
$ws.Range("M16").Value = -231.70587
$ws.Range("H22").Value = 768.8
$ws.Range("I22").Value = 601.6667
$ws.Range("K22").Value = 601.6667
$ws.Range("M22").Value = -306.6667
$ws.Range("H27").Value = 768.8
$ws.Range("I27").Value = 601.6667
$ws.Range("K27").Value = 601.6667
$ws.Range("M27").Value = -494.6667
$ws.Range("H40").Value = 1720.6757
$ws.Range("I40").Value = 1575.6562
$ws.Range("J40").Value = 2648.8
$ws.Range("K40").Value = 1575.6562
$ws.Range("L40").Value = 2648.8
$ws.Range("M40").Value = -1439.6562
$ws.Range("N40").Value = -2920.8
$ws.Range("H46").Value = 3591.682
$ws.Range("J46").Value = 3927.1
$ws.Range("L46").Value = 3927.1
$ws.Range("N46").Value = -4303.1
$ws.Range("H93").Value = 9348.75
$ws.Range("I93").Value = 9348.75
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 9348.75
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = None
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 4000
$ws.Range("I100").Value = 4000
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4000
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = None
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 387452.56
$ws.Range("I122").Value = 2680.8667
$ws.Range("J122").Value = 912141.25
$ws.Range("K122").Value = 8042.6001
$ws.Range("L122").Value = 2736423.75
$ws.Range("M122").Value = -5592.6001
$ws.Range("N122").Value = -2741323.75
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = None
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 2717.6667
$ws.Range("I132").Value = 1412.8889
$ws.Range("K132").Value = 4238.6667
$ws.Range("M132").Value = -1708.6667
$ws.Range("H136").Value = 335774.84
$ws.Range("I136").Value = 365499.8
$ws.Range("J136").Value = 8800.4
$ws.Range("K136").Value = 1096499.4
$ws.Range("L136").Value = 26401.2
$ws.Range("M136").Value = -1093949.4
$ws.Range("N136").Value = -31501.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 50041
$ws.Range("J44").Value = 50041
$ws.Range("L44").Value = 50041
$ws.Range("N44").Value = -51149
$ws.Range("H96").Value = 500999.5
$ws.Range("J96").Value = 2000
$ws.Range("L96").Value = 2000
$ws.Range("N96").Value = -4746
$ws.Range("H107").Value = 2233.6667
$ws.Range("I107").Value = 2233.6667
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 6701.000100000001
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = None
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 19610116
$ws.Range("I122").Value = 24391658
$ws.Range("J122").Value = 5793.8
$ws.Range("K122").Value = 73174974
$ws.Range("L122").Value = 17381.4
$ws.Range("M122").Value = -73172524
$ws.Range("N122").Value = -22281.4
$ws.Range("H126").Value = 1736.2222
$ws.Range("J126").Value = 3216.8
$ws.Range("L126").Value = 9650.400000000001
$ws.Range("N126").Value = -14590.4
$ws.Range("H132").Value = 21889.77
$ws.Range("I132").Value = 1141.4634
$ws.Range("J132").Value = 99224.37
$ws.Range("K132").Value = 3424.3902
$ws.Range("L132").Value = 297673.11
$ws.Range("M132").Value = -894.3902000000003
$ws.Range("N132").Value = -302733.11
$ws.Range("H136").Value = 267629.66
$ws.Range("I136").Value = 307187.8
$ws.Range("K136").Value = 921563.3999999999
$ws.Range("M136").Value = -919013.3999999999
